$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Testcase 3 (row 5): the "Expected" text was edited to drop the leading
# dash - "-Default Strategy Ethical selected" -> "Default Strategy Ethical selected"
$ws.Range("C5").Value = "Default Strategy Ethical selected"

# Reflect the cell the author had selected when the file was last saved.
$ws.Range("B13").Select() | Out-Null
